$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.133.72"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.905.91"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").Value = "253.73"
$ws.Range("E5").Value = "  +3.31%  "
$ws.Range("D6").Value = "0.697"
$ws.Range("E6").Value = "  +1.78%  "
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").Value = "41.49"
$ws.Range("E8").Value = "  +2.63%  "
$ws.Range("E9").Value = "  +4.29%  "
$ws.Range("D10").Value = "52.53"
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").Value = "0.0752"
$ws.Range("E11").Value = "  +4.97%  "
$ws.Range("D12").Value = "0.0980"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").Value = "13.20"
$ws.Range("E13").Value = "  +5.57%  "
$ws.Range("D14").Value = "2.181.68"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("E15").Value = "  +4.72%  "
$ws.Range("D16").Value = "5.03"
$ws.Range("E16").Value = "  +5.55%  "
$ws.Range("D17").Value = "1.902.58"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "35.120.24"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "73.78"
$ws.Range("E19").Value = "  +2.66%  "
$ws.Range("D20").Value = "0.0₃0839"
$ws.Range("E20").Value = "  +3.15%  "
$ws.Range("D21").Value = "243.05"
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("D22").Value = "12.99"
$ws.Range("E22").Value = "  +3.69%  "
$ws.Range("D23").Value = "5.05"
$ws.Range("E23").Value = "  +6.22%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").Value = "2.44"
$ws.Range("E25").Value = "  +5.35%  "
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").Value = "167.90"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "8.60"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").Value = "18.56"
$ws.Range("E29").Value = "  +2.13%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "4.128.42"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("E32").Value = "  +7.52%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "4.34"
$ws.Range("E33").Value = "  +5.02%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "2.01"
$ws.Range("E34").Value = "  +8.49%  "
$ws.Range("E35").Value = "  +7.71%  "
$ws.Range("D36").Value = "4.21"
$ws.Range("E36").Value = "  +3.38%  "
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D38").Value = "0.852"
$ws.Range("E38").Value = "  -5.50%  "
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("D40").Value = "103.41"
$ws.Range("E40").Value = "  +15.86%  "
$ws.Range("D41").Value = "17.26"
$ws.Range("E41").Value = "  +7.93%  "
$ws.Range("E43").Value = "  +1.72%  "
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").Value = "1.304.04"
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("D47").Value = "12.69"
$ws.Range("E47").Value = "  +4.24%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("D51").Value = "0.0751"
$ws.Range("E51").Value = "  +7.05%  "
